$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder header row (row 1): new column order by label
$ws.Range("A1").Value = "kitchens_1"
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("C1").Value = "living_rooms_1"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# Data rows 2-7: values follow their header label (columns reordered along with header)
$data = @(
    @(0,0,0,0,1,0),
    @(0,1,0,0,0,0),
    @(0,0,0,0,0,1),
    @(0,0,0,1,0,0),
    @(1,0,0,0,0,0),
    @(0,0,1,0,0,0)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rowValues = $data[$i]
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $col = $j + 1
        $ws.Cells.Item($row, $col).Value = $rowValues[$j]
    }
}
